$wb = $excel.ActiveWorkbook

$changesBySheet = @{}

$changesBySheet["ALC"] = @(
    @{Ref="H92"; Action="set"; Value=451},
    @{Ref="J92"; Action="set"; Value=805},
    @{Ref="L92"; Action="set"; Value=805},
    @{Ref="N92"; Action="set"; Value=-3301},
    @{Ref="H97"; Action="set"; Value=739.8},
    @{Ref="J97"; Action="set"; Value=739.8},
    @{Ref="L97"; Action="set"; Value=2219.4},
    @{Ref="N97"; Action="set"; Value=-3211.4},
    @{Ref="H98"; Action="set"; Value=660.36365},
    @{Ref="I98"; Action="set"; Value=621.9},
    @{Ref="J98"; Action="set"; Value=1045},
    @{Ref="K98"; Action="set"; Value=621.9},
    @{Ref="L98"; Action="set"; Value=1045},
    @{Ref="M98"; Action="set"; Value=876.1},
    @{Ref="N98"; Action="set"; Value=-4041},
    @{Ref="H101"; Action="set"; Value=3911.6667},
    @{Ref="I101"; Action="set"; Value=3911.6667},
    @{Ref="K101"; Action="set"; Value=11735.0001},
    @{Ref="M101"; Action="set"; Value=-10113.0001},
    @{Ref="H113"; Action="set"; Value=2947.5},
    @{Ref="I113"; Action="set"; Value=0},
    @{Ref="K113"; Action="set"; Value=0},
    @{Ref="M113"; Action="clear"; Value=$null},
    @{Ref="H116"; Action="set"; Value=3499.5},
    @{Ref="I116"; Action="set"; Value=3000},
    @{Ref="K116"; Action="set"; Value=3000},
    @{Ref="M116"; Action="set"; Value=442},
    @{Ref="H122"; Action="set"; Value=660.36365},
    @{Ref="I122"; Action="set"; Value=621.9},
    @{Ref="J122"; Action="set"; Value=1045},
    @{Ref="K122"; Action="set"; Value=1865.7},
    @{Ref="L122"; Action="set"; Value=3135},
    @{Ref="M122"; Action="set"; Value=584.3000000000002},
    @{Ref="N122"; Action="set"; Value=-8035},
    @{Ref="H137"; Action="set"; Value=1625.2609},
    @{Ref="I137"; Action="set"; Value=1493.3889},
    @{Ref="K137"; Action="set"; Value=4480.1667},
    @{Ref="M137"; Action="set"; Value=-1930.1667}
)

$changesBySheet["ARM"] = @(
    @{Ref="H61"; Action="set"; Value=9574.5},
    @{Ref="I61"; Action="set"; Value=9574.5},
    @{Ref="K61"; Action="set"; Value=9574.5},
    @{Ref="M61"; Action="set"; Value=-9362.5},
    @{Ref="H74"; Action="set"; Value=2258},
    @{Ref="I74"; Action="set"; Value=1011},
    @{Ref="J74"; Action="set"; Value=5999},
    @{Ref="K74"; Action="set"; Value=1011},
    @{Ref="L74"; Action="set"; Value=5999},
    @{Ref="M74"; Action="set"; Value=-137},
    @{Ref="N74"; Action="set"; Value=-7747},
    @{Ref="H77"; Action="set"; Value=2258},
    @{Ref="I77"; Action="set"; Value=1011},
    @{Ref="J77"; Action="set"; Value=5999},
    @{Ref="K77"; Action="set"; Value=5055},
    @{Ref="L77"; Action="set"; Value=29995},
    @{Ref="M77"; Action="set"; Value=-687},
    @{Ref="N77"; Action="set"; Value=-38731},
    @{Ref="H88"; Action="set"; Value=5344.3335},
    @{Ref="I88"; Action="set"; Value=4251},
    @{Ref="J88"; Action="set"; Value=6219},
    @{Ref="K88"; Action="set"; Value=4251},
    @{Ref="L88"; Action="set"; Value=6219},
    @{Ref="M88"; Action="set"; Value=-3845},
    @{Ref="N88"; Action="set"; Value=-7031},
    @{Ref="H91"; Action="set"; Value=5344.3335},
    @{Ref="I91"; Action="set"; Value=4251},
    @{Ref="J91"; Action="set"; Value=6219},
    @{Ref="K91"; Action="set"; Value=4251},
    @{Ref="L91"; Action="set"; Value=6219},
    @{Ref="M91"; Action="set"; Value=-2847},
    @{Ref="N91"; Action="set"; Value=-9027},
    @{Ref="H136"; Action="set"; Value=9574.5},
    @{Ref="I136"; Action="set"; Value=9574.5},
    @{Ref="K136"; Action="set"; Value=28723.5},
    @{Ref="M136"; Action="set"; Value=-26173.5}
)

$changesBySheet["BSM"] = @(
    @{Ref="H54"; Action="set"; Value=5126.7144},
    @{Ref="I54"; Action="set"; Value=5126.7144},
    @{Ref="K54"; Action="set"; Value=5126.7144},
    @{Ref="M54"; Action="set"; Value=-4642.7144}
)

$changesBySheet["CRP"] = @(
    @{Ref="H16"; Action="set"; Value=1033.3334},
    @{Ref="J16"; Action="set"; Value=1100},
    @{Ref="L16"; Action="set"; Value=1100},
    @{Ref="N16"; Action="set"; Value=-1674},
    @{Ref="H31"; Action="set"; Value=4903.727},
    @{Ref="J31"; Action="set"; Value=9004},
    @{Ref="L31"; Action="set"; Value=9004},
    @{Ref="N31"; Action="set"; Value=-9594},
    @{Ref="H34"; Action="set"; Value=4903.727},
    @{Ref="J34"; Action="set"; Value=9004},
    @{Ref="L34"; Action="set"; Value=9004},
    @{Ref="N34"; Action="set"; Value=-9408},
    @{Ref="H58"; Action="set"; Value=5090.273},
    @{Ref="I58"; Action="set"; Value=3599.3},
    @{Ref="K58"; Action="set"; Value=3599.3},
    @{Ref="M58"; Action="set"; Value=-3396.3},
    @{Ref="H99"; Action="set"; Value=5531.385},
    @{Ref="I99"; Action="set"; Value=4898.909},
    @{Ref="J99"; Action="set"; Value=9010},
    @{Ref="K99"; Action="set"; Value=4898.909},
    @{Ref="L99"; Action="set"; Value=9010},
    @{Ref="M99"; Action="set"; Value=-3400.909},
    @{Ref="N99"; Action="set"; Value=-12006},
    @{Ref="H113"; Action="set"; Value=1033.3334},
    @{Ref="J113"; Action="set"; Value=1100},
    @{Ref="L113"; Action="set"; Value=1100},
    @{Ref="N113"; Action="set"; Value=-5440},
    @{Ref="H126"; Action="set"; Value=5531.385},
    @{Ref="I126"; Action="set"; Value=4898.909},
    @{Ref="J126"; Action="set"; Value=9010},
    @{Ref="K126"; Action="set"; Value=14696.727},
    @{Ref="L126"; Action="set"; Value=27030},
    @{Ref="M126"; Action="set"; Value=-12226.727},
    @{Ref="N126"; Action="set"; Value=-31970},
    @{Ref="H136"; Action="set"; Value=5090.273},
    @{Ref="I136"; Action="set"; Value=3599.3},
    @{Ref="K136"; Action="set"; Value=10797.9},
    @{Ref="M136"; Action="set"; Value=-8247.900000000001}
)

$changesBySheet["CUL"] = @(
    @{Ref="H2"; Action="set"; Value=120.5},
    @{Ref="J2"; Action="set"; Value=120.5},
    @{Ref="L2"; Action="set"; Value=723},
    @{Ref="N2"; Action="set"; Value=-949},
    @{Ref="H122"; Action="set"; Value=794.25},
    @{Ref="I122"; Action="set"; Value=350},
    @{Ref="J122"; Action="set"; Value=1238.5},
    @{Ref="K122"; Action="set"; Value=3150},
    @{Ref="L122"; Action="set"; Value=11146.5},
    @{Ref="M122"; Action="set"; Value=-700},
    @{Ref="N122"; Action="set"; Value=-16046.5},
    @{Ref="H125"; Action="set"; Value=3000},
    @{Ref="J125"; Action="set"; Value=3000},
    @{Ref="L125"; Action="set"; Value=9000},
    @{Ref="N125"; Action="set"; Value=-18840},
    @{Ref="H131"; Action="set"; Value=907.94116},
    @{Ref="J131"; Action="set"; Value=931.29034},
    @{Ref="L131"; Action="set"; Value=2793.87102},
    @{Ref="N131"; Action="set"; Value=-12873.87102},
    @{Ref="H137"; Action="set"; Value=1997},
    @{Ref="I137"; Action="set"; Value=1997},
    @{Ref="J137"; Action="set"; Value=0},
    @{Ref="K137"; Action="set"; Value=5991},
    @{Ref="L137"; Action="set"; Value=0},
    @{Ref="M137"; Action="set"; Value=-891},
    @{Ref="N137"; Action="clear"; Value=$null}
)

$changesBySheet["GSM"] = @(
    @{Ref="H4"; Action="set"; Value=1000},
    @{Ref="J4"; Action="set"; Value=1000},
    @{Ref="L4"; Action="set"; Value=1000},
    @{Ref="N4"; Action="set"; Value=-1224},
    @{Ref="H132"; Action="set"; Value=4423.3335},
    @{Ref="I132"; Action="set"; Value=2262.8},
    @{Ref="K132"; Action="set"; Value=6788.400000000001},
    @{Ref="M132"; Action="set"; Value=-4258.400000000001}
)

$changesBySheet["LTW"] = @(
    @{Ref="H7"; Action="set"; Value=5996.25},
    @{Ref="I7"; Action="set"; Value=4995},
    @{Ref="K7"; Action="set"; Value=4995},
    @{Ref="M7"; Action="set"; Value=-4883},
    @{Ref="H111"; Action="set"; Value=22693.5},
    @{Ref="J111"; Action="set"; Value=22693.5},
    @{Ref="L111"; Action="set"; Value=22693.5},
    @{Ref="N111"; Action="set"; Value=-30873.5},
    @{Ref="H126"; Action="set"; Value=5996.25},
    @{Ref="I126"; Action="set"; Value=4995},
    @{Ref="K126"; Action="set"; Value=14985},
    @{Ref="M126"; Action="set"; Value=-12515}
)

foreach ($sheetName in $changesBySheet.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($c in $changesBySheet[$sheetName]) {
        if ($c.Action -eq "clear") {
            $ws.Range($c.Ref).ClearContents()
        } else {
            $ws.Range($c.Ref).Value = $c.Value
        }
    }
}

Write-Host "Applied cell changes across $($changesBySheet.Keys.Count) sheets"
